# Generate Report for Handoff
# Adds a new handed-off file (3d1c2947-6b92-45fc-9532-6552f3cbe3ce.md) as a
# new row to each of the three tracking sheets: Overview, zh-cn, de-de.
#
# Every text value is written through `.Formula = "'" + value` (Excel's
# "force text" apostrophe prefix) so that literal strings which otherwise
# look like booleans/dates/numbers ("True", "False", "2016-08-26 14:50:51",
# "0001-01-01 00:00:00", or even the empty string) are stored as plain text
# instead of being auto-coerced to a typed value.

function Set-Text($cell, [string]$text) {
    $cell.Formula = "'" + $text
}

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39e579a4c5a37cdc16c310277159f37ec5fd4cfd/e2e/"
$newFile  = "3d1c2947-6b92-45fc-9532-6552f3cbe3ce.md"
$newPath  = "e2e\" + $newFile

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | Path And Name | Extension | Publish URL |
#                   zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$rowOverview = $loOverview.ListRows.Add()

Set-Text $rowOverview.Range.Cells.Item(1, 1) $newFile
$bCell = $rowOverview.Range.Cells.Item(1, 2)
Set-Text $bCell $newPath
$wsOverview.Hyperlinks.Add($bCell, ($repoBase + $newFile), "", "", $newPath)
Set-Text $rowOverview.Range.Cells.Item(1, 3) ".md"
Set-Text $rowOverview.Range.Cells.Item(1, 4) ""
Set-Text $rowOverview.Range.Cells.Item(1, 5) "Ready for handoff"
Set-Text $rowOverview.Range.Cells.Item(1, 6) "Ready for handoff"
Set-Text $rowOverview.Range.Cells.Item(1, 7) "2016-08-26 14:50:51"

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Source Path |
#   Priority | Content Duplicate | Latest Handoff File | Latest Handoff
#   Datetime | Latest Target File | Latest Handback File | Latest Handback
#   DateTime | Reference Tokens | To be localized | Dependency From |
#   Has metadata | Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$rowZhCn = $loZhCn.ListRows.Add()

$aCellZh = $rowZhCn.Range.Cells.Item(1, 1)
Set-Text $aCellZh $newFile
$wsZhCn.Hyperlinks.Add($aCellZh, ($repoBase + $newFile), "", "", $newFile)
Set-Text $rowZhCn.Range.Cells.Item(1, 2)  ".md"
Set-Text $rowZhCn.Range.Cells.Item(1, 3)  "Ready for handoff"
Set-Text $rowZhCn.Range.Cells.Item(1, 4)  "e2e"
Set-Text $rowZhCn.Range.Cells.Item(1, 5)  "ht"
Set-Text $rowZhCn.Range.Cells.Item(1, 6)  "False"
Set-Text $rowZhCn.Range.Cells.Item(1, 7)  "3d1c2947-6b92-45fc-9532-6552f3cbe3ce.f67035688b80528b8e27d5f28006c8008ae1e32b.zh-cn.xlf"
Set-Text $rowZhCn.Range.Cells.Item(1, 8)  "2016-08-26 14:50:46"
Set-Text $rowZhCn.Range.Cells.Item(1, 9)  ""
Set-Text $rowZhCn.Range.Cells.Item(1, 10) ""
Set-Text $rowZhCn.Range.Cells.Item(1, 11) "0001-01-01 00:00:00"
Set-Text $rowZhCn.Range.Cells.Item(1, 12) ""
Set-Text $rowZhCn.Range.Cells.Item(1, 13) "True"
Set-Text $rowZhCn.Range.Cells.Item(1, 14) ""
Set-Text $rowZhCn.Range.Cells.Item(1, 15) "False"
Set-Text $rowZhCn.Range.Cells.Item(1, 16) ""

# ---------------------------------------------------------------------
# Sheet "de-de": same column layout as zh-cn
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$rowDeDe = $loDeDe.ListRows.Add()

$aCellDe = $rowDeDe.Range.Cells.Item(1, 1)
Set-Text $aCellDe $newFile
$wsDeDe.Hyperlinks.Add($aCellDe, ($repoBase + $newFile), "", "", $newFile)
Set-Text $rowDeDe.Range.Cells.Item(1, 2)  ".md"
Set-Text $rowDeDe.Range.Cells.Item(1, 3)  "Ready for handoff"
Set-Text $rowDeDe.Range.Cells.Item(1, 4)  "e2e"
Set-Text $rowDeDe.Range.Cells.Item(1, 5)  "ht"
Set-Text $rowDeDe.Range.Cells.Item(1, 6)  "False"
Set-Text $rowDeDe.Range.Cells.Item(1, 7)  "3d1c2947-6b92-45fc-9532-6552f3cbe3ce.f67035688b80528b8e27d5f28006c8008ae1e32b.de-de.xlf"
Set-Text $rowDeDe.Range.Cells.Item(1, 8)  "2016-08-26 14:50:51"
Set-Text $rowDeDe.Range.Cells.Item(1, 9)  ""
Set-Text $rowDeDe.Range.Cells.Item(1, 10) ""
Set-Text $rowDeDe.Range.Cells.Item(1, 11) "0001-01-01 00:00:00"
Set-Text $rowDeDe.Range.Cells.Item(1, 12) ""
Set-Text $rowDeDe.Range.Cells.Item(1, 13) "True"
Set-Text $rowDeDe.Range.Cells.Item(1, 14) ""
Set-Text $rowDeDe.Range.Cells.Item(1, 15) "False"
Set-Text $rowDeDe.Range.Cells.Item(1, 16) ""

Write-Host "Added handoff row for" $newFile "to Overview, zh-cn, de-de"
